$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "MODEL_CONDITION" header text to "MODELCONDITION" (shared string text change)
$ws.Cells.Item(1, 5).Value = "MODELCONDITION"

# Delete column A entirely, shifting B:F left to A:E
$ws.Columns("A").Delete()
